$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "- Reologia de fluidos,- Dimensionamento de tubulações,- Acessórios e bombeamento para fluidos industriais- Agitação e mistura- Caracterização de partículas e leito de partículas- Sedimentação- Filtração- Processos com membranas- Operações unitárias de troca térmica: trocadores de calor e evaporadores- Psicrometria: carta psicrométrica e propriedades do ar",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- Reologia de fluidos,^l- Dimensionamento de tubulações,^l- Acessórios e bombeamento para fluidos industriais^l- Agitação e mistura^l- Caracterização de partículas e leito de partículas^l- Sedimentação^l- Filtração^l- Processos com membranas^l- Operações unitárias de troca térmica: trocadores de calor e evaporadores^l- Psicrometria: carta psicrométrica e propriedades do ar",
    2
)

$d.Content.Find.Execute(
    "- Fluid rheology- Sizing of pipes- Accessories and pumping for industrial fluids- Stirring and mixing- Particle characterization and particle bed- Sedimentation- Filtration- Processes with membranes- Unit heat exchange operations: heat exchangers and evaporators- Psychrometry: psychrometric chart and air properties",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- Fluid rheology^l- Sizing of pipes^l- Accessories and pumping for industrial fluids^l- Stirring and mixing^l- Particle characterization and particle bed^l- Sedimentation^l- Filtration^l- Processes with membranes^l- Unit heat exchange operations: heat exchangers and evaporators^l- Psychrometry: psychrometric chart and air properties",
    2
)

$d.Content.Find.Execute(
    "FOUST, A.S., WENZEL, L. A., CLUMP, C.W., MAUS, L., ANDERSEN, L.B. Princípio das operações unitárias. Rio de Janeiro: Editora Guanabara Dois, 1982.GEANKOPLIS, C.J. Procesos de transporte y operaciones unitarias. Compañía Editorial Continental, S.A. de C.V. México, D.F., 1998.PERRY, R.H. and CHILTON, C.H. Manual de Engenharia Química. 5a ed., Guanabara Dois, Rio de Janeiro, 1986.REYNOLDS, T.D.; RICHARDS, P. Unit Operations and Processes in environmental Engineering. PWS Publishing, 1996.MACINTYRE, A.J. Bombas e Instalações de Bombeamento. LTC, Rio de Janeiro, 1997",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "FOUST, A.S., WENZEL, L. A., CLUMP, C.W., MAUS, L., ANDERSEN, L.B. Princípio das operações unitárias. Rio de Janeiro: Editora Guanabara Dois, 1982.^lGEANKOPLIS, C.J. Procesos de transporte y operaciones unitarias. Compañía Editorial Continental, S.A. de C.V. México, D.F., 1998.^lPERRY, R.H. and CHILTON, C.H. Manual de Engenharia Química. 5a ed., Guanabara Dois, Rio de Janeiro, 1986.^lREYNOLDS, T.D.; RICHARDS, P. Unit Operations and Processes in environmental Engineering. PWS Publishing, 1996.^lMACINTYRE, A.J. Bombas e Instalações de Bombeamento. LTC, Rio de Janeiro, 1997",
    2
)
